$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.680.82'
$ws.Range("E2").Value = '  +5.19%  '
$ws.Range("D3").Value = '2.269.26'
$ws.Range("E3").Value = '  +2.98%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.68'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +1.60%  '
$ws.Range("E6").Value = '  +2.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.56'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +6.77%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.430'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +7.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.103'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +16.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.45'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.95'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +16.93%  '
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").Value = '2.607.79'
$ws.Range("E14").Value = '  +2.83%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.70'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +2.47%  '
$ws.Range("E16").Value = '  +5.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.826'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +4.45%  '
$ws.Range("D18").Value = '2.272.78'
$ws.Range("E18").Value = '  +2.94%  '
$ws.Range("D19").Value = '43.520.67'
$ws.Range("D20").Value = '0.0₃0999'
$ws.Range("E20").Value = '  +11.32%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.75'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +2.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.10'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +1.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '249.82'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +3.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.50'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +7.15%  '
$ws.Range("E26").Value = '  +1.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.89'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +2.87%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '172.32'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +2.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '21.00'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +6.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.138'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  -1.01%  '
$ws.Range("E31").Value = '  +1.62%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.80'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +10.43%  '
$ws.Range("E33").Value = '  +1.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0686'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +6.46%  '
$ws.Range("E35").Value = '  +2.57%  '
$ws.Range("E36").Value = '  +3.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.81'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  +5.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.82'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +7.53%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.33'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  -0.91%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0249'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +5.95%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.42'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  -1.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.34'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +6.36%  '
$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.45'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +4.18%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0960'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +0.86%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.21'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +0.95%  '
$ws.Range("B47").Value = 'Celestia'
$ws.Range("C47").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.35'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +21.44%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '97.63'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  +0.98%  '
$ws.Range("D49").Value = '1.477.39'
$ws.Range("E49").Value = '  +1.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.35'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +5.80%  '
$ws.Range("E51").Value = '  +1.16%  '
Write-Host "Applied cryptos update"
